$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 12 (pushes old row 12 down to row 13),
# inheriting formatting from the row above it (row 11).
$ws.Rows.Item(12).Insert()

# Update C11 value
$ws.Range("C11").Value = 84.212092567146101

# Fill the newly inserted row 12, copying the style of row 11 (the row
# above) so it keeps the normal (non-thick-bottom) borders.
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122) | Out-Null

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Decision Tree Regression+feature selection"
$ws.Range("C12").Value = 91.0827149532699

# Update old row 12 (now row 13) - the Id number should now be 11
$ws.Range("A13").Value = 11

# Update the selection to match target state
$ws.Range("C17").Select()
